$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value from 10.0 to 1.0
$ws.Range("A2").Value = 1.0

# Update F2 timestamp string
$ws.Range("F2").Value = "Sat May 16 13:51:18 MSK 2020"

# Delete row 3 entirely (was A3:F3)
$ws.Rows.Item(3).Delete()
